$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wally_BMS")

# Fill in the newly added BOM line (row 40) - connector socket that was
# previously blank. Values are entered in this particular order so the
# new shared-string table entries come out in the same order Excel
# produced them in.
$ws.Range("D40").Value = "900-1727181111CT-ND"
$ws.Range("B40").Value = "J2_2"
$ws.Range("A40").Value = "CONN SOCKET 16AWG CRIMP TIN"
$ws.Range("C40").Value = 10
$ws.Range("E40").Value = 1727181111
$ws.Range("F40").Value = "Molex"

# Setting .Value above cleared the "quote prefix" formatting that A40/B40
# (text cells with a wrapping border) originally carried. Restore their
# original cell format (borrowed from the identically formatted row
# above) without touching the values we just entered.
$ws.Range("A39").Copy()
$ws.Range("A40").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B39").Copy()
$ws.Range("B40").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the view state to match where the user scrolled/selected after
# the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A40").Select()
